$d = $word.ActiveDocument

# Locate the very end of the document body content (just before sectPr)
$endPos = $d.Content.End
$r = $d.Range($endPos, $endPos)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:pPr>' +
       '<w:ind w:left="360" w:firstLine="0"/>' +
       '<w:rPr><w:b/><w:bCs/></w:rPr>' +
       '</w:pPr>' +
       '<w:r>' +
       '<w:rPr><w:b/><w:bCs/></w:rPr>' +
       '<w:t>https://github.com/Robertino2809/Programsko-in-enjerstvo---projekt</w:t>' +
       '</w:r>' +
       '</w:p>'

$r.InsertXML($xml)

# InsertXML silently drops w:firstLine="0" since it is the implicit default;
# force it back onto the new paragraph explicitly so it round-trips in the OOXML.
$newCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newCount)
$newPara.Format.FirstLineIndent = 0
